$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(42, 1).Value = "TZP"
$ws.Cells.Item(42, 2).Value = "Year"
$ws.Cells.Item(42, 3).Value = [double]"1.367286867747675e-25"

$ws.Cells.Item(43, 1).Value = "TZP"
$ws.Cells.Item(43, 2).Value = "Specimen_type"
$ws.Cells.Item(43, 3).Value = [double]"7.481013117150814e-07"

$ws.Cells.Item(44, 1).Value = "TZP"
$ws.Cells.Item(44, 2).Value = "Gender"
$ws.Cells.Item(44, 3).Value = [double]"0.003250545772315636"

$ws.Cells.Item(45, 1).Value = "TZP"
$ws.Cells.Item(45, 2).Value = "Hospital:Ward_ED_ICU"
$ws.Cells.Item(45, 3).Value = [double]"6.377440382755414e-10"
